# "update to r4 qa"
# - bump CapabilityStatement version / fhirVersion / IG canonical URL on the
#   "meta" sheet from the STU3 (1.0.0 / FHIR 3.0.1) coordinates to the R4
#   (1.1.0 / FHIR 4.0.0) coordinates
# - drop the "/STU3" path segment from every Profile canonical URL on the
#   "profiles" sheet, since the profiles now live directly under
#   .../davinci-deqm/StructureDefinition/...

$wb = $excel.ActiveWorkbook

# --- meta sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("meta")
$meta.Range("B2").Value = "1.1.0"
$meta.Range("B3").Value = "4.0.0"
$meta.Range("B6").Value = "http://hl7.org/fhir/us/davinci-deqm/ImplementationGuide/hl7.fhir.us.davinci-deqm-1.1.0"

# --- profiles sheet -----------------------------------------------------
$profiles = $wb.Worksheets.Item("profiles")
for ($r = 2; $r -le 9; $r++) {
    $cell = $profiles.Cells.Item($r, 1)
    $old = [string]$cell.Value2
    $cell.Value = $old.Replace("/STU3/StructureDefinition/", "/StructureDefinition/")
}

# --- view state: active sheet/selection moved from "ops" to "meta" ------
$ops = $wb.Worksheets.Item("ops")
$ops.Activate()
$ops.Range("B7").Select()

$profiles.Activate()
$profiles.Range("A16").Select()

$meta.Activate()
$meta.Range("B7").Select()
